# Auto-generated Excel COM-interop script applying the Goblin_Profits.xlsx diff.
# Updates numeric leve-profit cells (columns H-N) across all 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the target workbook state.

$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2260
$ws.Range("I86").Value = 3100
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 3100
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -1977
$ws.Range("N86").Value = -3946
$ws.Range("H88").Value = 5034.6
$ws.Range("J88").Value = 5684.8823
$ws.Range("L88").Value = 5684.8823
$ws.Range("N88").Value = -6496.8823
$ws.Range("H89").Value = 2260
$ws.Range("I89").Value = 3100
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 15500
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -9884
$ws.Range("N89").Value = -19732
$ws.Range("H91").Value = 5034.6
$ws.Range("J91").Value = 5684.8823
$ws.Range("L91").Value = 5684.8823
$ws.Range("N91").Value = -8492.882300000001
$ws.Range("H92").Value = 1170.9524
$ws.Range("I92").Value = 975.94116
$ws.Range("K92").Value = 975.94116
$ws.Range("M92").Value = 272.05884
$ws.Range("H112").Value = 1472.0625
$ws.Range("J112").Value = 1546.5
$ws.Range("L112").Value = 4639.5
$ws.Range("N112").Value = -6855.5
$ws.Range("H116").Value = 6507.857
$ws.Range("I116").Value = 5542.625
$ws.Range("J116").Value = 7794.8335
$ws.Range("K116").Value = 5542.625
$ws.Range("L116").Value = 7794.8335
$ws.Range("M116").Value = -2100.625
$ws.Range("N116").Value = -14678.8335
$ws.Range("H127").Value = 1916.6666
$ws.Range("I127").Value = 1282.375
$ws.Range("J127").Value = 2641.5715
$ws.Range("K127").Value = 3847.125
$ws.Range("L127").Value = 7924.7145
$ws.Range("M127").Value = 1112.875
$ws.Range("N127").Value = -17844.7145
$ws.Range("H129").Value = 1189.0769
$ws.Range("I129").Value = 862
$ws.Range("K129").Value = 2586
$ws.Range("M129").Value = 2414
$ws.Range("H132").Value = 2383137.5
$ws.Range("J132").Value = 12502210
$ws.Range("L132").Value = 37506630
$ws.Range("N132").Value = -37511690
$ws.Range("H133").Value = 121998.336
$ws.Range("J133").Value = 121998.336
$ws.Range("L133").Value = 121998.336
$ws.Range("N133").Value = -132118.336
$ws.Range("H136").Value = 338957
$ws.Range("J136").Value = 338957
$ws.Range("L136").Value = 338957
$ws.Range("N136").Value = -349157
$ws.Range("H137").Value = 1088245.8
$ws.Range("I137").Value = 1138.4546
$ws.Range("J137").Value = 2416932.5
$ws.Range("K137").Value = 3415.3638
$ws.Range("L137").Value = 7250797.5
$ws.Range("M137").Value = -865.3638000000001
$ws.Range("N137").Value = -7255897.5

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 183767.47
$ws.Range("I32").Value = 183767.47
$ws.Range("K32").Value = 183767.47
$ws.Range("M32").Value = -183480.47
$ws.Range("H74").Value = 2514.7368
$ws.Range("I74").Value = 2099.3462
$ws.Range("K74").Value = 2099.3462
$ws.Range("M74").Value = -1225.3462
$ws.Range("H77").Value = 2514.7368
$ws.Range("I77").Value = 2099.3462
$ws.Range("K77").Value = 10496.731
$ws.Range("M77").Value = -6128.731
$ws.Range("H102").Value = 3394.7942
$ws.Range("I102").Value = 1286.1
$ws.Range("K102").Value = 1286.1
$ws.Range("M102").Value = 335.9000000000001
$ws.Range("H133").Value = 44495.75
$ws.Range("J133").Value = 44495.75
$ws.Range("L133").Value = 44495.75
$ws.Range("N133").Value = -49555.75

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 63331.332
$ws.Range("I2").Value = 54000
$ws.Range("J2").Value = 67997
$ws.Range("K2").Value = 54000
$ws.Range("L2").Value = 67997
$ws.Range("M2").Value = -53887
$ws.Range("N2").Value = -68223
$ws.Range("H75").Value = 33842.2
$ws.Range("I75").Value = 19607
$ws.Range("J75").Value = 43332.332
$ws.Range("K75").Value = 19607
$ws.Range("L75").Value = 43332.332
$ws.Range("M75").Value = -18671
$ws.Range("N75").Value = -45204.332
$ws.Range("H78").Value = 33842.2
$ws.Range("I78").Value = 19607
$ws.Range("J78").Value = 43332.332
$ws.Range("K78").Value = 58821
$ws.Range("L78").Value = 129996.996
$ws.Range("M78").Value = -54141
$ws.Range("N78").Value = -139356.996
$ws.Range("H82").Value = 100031910
$ws.Range("J82").Value = 250047500
$ws.Range("L82").Value = 250047500
$ws.Range("N82").Value = -250048266
$ws.Range("H85").Value = 100031910
$ws.Range("J85").Value = 250047500
$ws.Range("L85").Value = 250047500
$ws.Range("N85").Value = -250050152
$ws.Range("H86").Value = 11905588
$ws.Range("I86").Value = 886.03125
$ws.Range("J86").Value = 50000636
$ws.Range("K86").Value = 886.03125
$ws.Range("L86").Value = 50000636
$ws.Range("M86").Value = 236.96875
$ws.Range("N86").Value = -50002882
$ws.Range("H89").Value = 11905588
$ws.Range("I89").Value = 886.03125
$ws.Range("J89").Value = 50000636
$ws.Range("K89").Value = 4430.15625
$ws.Range("L89").Value = 250003180
$ws.Range("M89").Value = 1185.84375
$ws.Range("N89").Value = -250014412
$ws.Range("H94").Value = 1662.5264
$ws.Range("I94").Value = 1411.0883
$ws.Range("K94").Value = 1411.0883
$ws.Range("M94").Value = -960.0882999999999
$ws.Range("H96").Value = 23894.4
$ws.Range("I96").Value = 17105.445
$ws.Range("J96").Value = 84995
$ws.Range("K96").Value = 17105.445
$ws.Range("L96").Value = 84995
$ws.Range("M96").Value = -14359.445
$ws.Range("N96").Value = -90487
$ws.Range("H97").Value = 15799.6
$ws.Range("J97").Value = 17999.666
$ws.Range("L97").Value = 17999.666
$ws.Range("N97").Value = -19981.666
$ws.Range("H101").Value = 49999
$ws.Range("J101").Value = 49999
$ws.Range("L101").Value = 49999
$ws.Range("N101").Value = -56489
$ws.Range("H105").Value = 2638.5557
$ws.Range("I105").Value = 1949.6
$ws.Range("J105").Value = 3499.75
$ws.Range("K105").Value = 1949.6
$ws.Range("L105").Value = 3499.75
$ws.Range("M105").Value = -202.5999999999999
$ws.Range("N105").Value = -6993.75
$ws.Range("H132").Value = 197463.88
$ws.Range("J132").Value = 197463.88
$ws.Range("L132").Value = 197463.88
$ws.Range("N132").Value = -207583.88

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 7689.75
$ws.Range("J25").Value = 14749.5
$ws.Range("L25").Value = 14749.5
$ws.Range("N25").Value = -15097.5
$ws.Range("H31").Value = 3315.3333
$ws.Range("I31").Value = 1280.3334
$ws.Range("J31").Value = 4332.8335
$ws.Range("K31").Value = 1280.3334
$ws.Range("L31").Value = 4332.8335
$ws.Range("M31").Value = -985.3334
$ws.Range("N31").Value = -4922.8335
$ws.Range("H34").Value = 3315.3333
$ws.Range("I34").Value = 1280.3334
$ws.Range("J34").Value = 4332.8335
$ws.Range("K34").Value = 1280.3334
$ws.Range("L34").Value = 4332.8335
$ws.Range("M34").Value = -1078.3334
$ws.Range("N34").Value = -4736.8335
$ws.Range("H105").Value = 6091.273
$ws.Range("I105").Value = 3820.4285
$ws.Range("K105").Value = 3820.4285
$ws.Range("M105").Value = -2073.4285
$ws.Range("H132").Value = 1637.975
$ws.Range("I132").Value = 974.8857400000001
$ws.Range("J132").Value = 6279.6
$ws.Range("K132").Value = 2924.65722
$ws.Range("L132").Value = 18838.8
$ws.Range("M132").Value = -394.6572200000001
$ws.Range("N132").Value = -23898.8
$ws.Range("H135").Value = 121821.664
$ws.Range("J135").Value = 121821.664
$ws.Range("L135").Value = 121821.664
$ws.Range("N135").Value = -131961.664

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 179.94737
$ws.Range("I12").Value = 39.333332
$ws.Range("K12").Value = 117.999996
$ws.Range("M12").Value = 55.000004
$ws.Range("H23").Value = 2128.3845
$ws.Range("I23").Value = 10036
$ws.Range("J23").Value = 690.63635
$ws.Range("K23").Value = 30108
$ws.Range("L23").Value = 2071.90905
$ws.Range("M23").Value = -29873
$ws.Range("N23").Value = -2541.90905
$ws.Range("H34").Value = 3446.1462
$ws.Range("J34").Value = 3446.1462
$ws.Range("L34").Value = 10338.4386
$ws.Range("N34").Value = -10506.4386
$ws.Range("H39").Value = 3043.818
$ws.Range("J39").Value = 2370.4285
$ws.Range("L39").Value = 7111.2855
$ws.Range("N39").Value = -7699.2855
$ws.Range("H40").Value = 93
$ws.Range("I40").Value = 84.09999999999999
$ws.Range("J40").Value = 137.5
$ws.Range("K40").Value = 336.4
$ws.Range("L40").Value = 550
$ws.Range("M40").Value = -267.4
$ws.Range("N40").Value = -688
$ws.Range("H51").Value = 893.5714
$ws.Range("I51").Value = 701
$ws.Range("K51").Value = 2103
$ws.Range("M51").Value = -1643
$ws.Range("H55").Value = 335041.88
$ws.Range("I55").Value = 833653.2
$ws.Range("J55").Value = 2634.3333
$ws.Range("K55").Value = 2500959.6
$ws.Range("L55").Value = 7902.999899999999
$ws.Range("M55").Value = -2500782.6
$ws.Range("N55").Value = -8256.999899999999
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value = 500050000
$ws.Range("I97").Value = 1000000000
$ws.Range("K97").Value = 3000000000
$ws.Range("M97").Value = -2999999504
$ws.Range("H134").Value = 5678.174
$ws.Range("I134").Value = 1674.5
$ws.Range("J134").Value = 7091.2354
$ws.Range("K134").Value = 5023.5
$ws.Range("L134").Value = 21273.7062
$ws.Range("M134").Value = 46.5
$ws.Range("N134").Value = -31413.7062

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7666.3335
$ws.Range("J7").Value = 7666.3335
$ws.Range("L7").Value = 7666.3335
$ws.Range("N7").Value = -7890.3335
$ws.Range("H8").Value = 7666.3335
$ws.Range("J8").Value = 7666.3335
$ws.Range("L8").Value = 7666.3335
$ws.Range("N8").Value = -7944.3335
$ws.Range("H10").Value = 257333.33
$ws.Range("J10").Value = 257333.33
$ws.Range("L10").Value = 257333.33
$ws.Range("N10").Value = -257671.33
$ws.Range("H12").Value = 383908.5
$ws.Range("I12").Value = 8221.5
$ws.Range("K12").Value = 8221.5
$ws.Range("M12").Value = -8081.5
$ws.Range("H14").Value = 58875.668
$ws.Range("I14").Value = 102446.2
$ws.Range("J14").Value = 4412.5
$ws.Range("K14").Value = 102446.2
$ws.Range("L14").Value = 4412.5
$ws.Range("M14").Value = -102278.2
$ws.Range("N14").Value = -4748.5
$ws.Range("H36").Value = 20000
$ws.Range("J36").Value = 20000
$ws.Range("L36").Value = 20000
$ws.Range("N36").Value = -20970
$ws.Range("H43").Value = 2928.25
$ws.Range("I43").Value = 2928.25
$ws.Range("K43").Value = 2928.25
$ws.Range("M43").Value = -2777.25
$ws.Range("H80").Value = 41670496
$ws.Range("J80").Value = 6492.143
$ws.Range("L80").Value = 6492.143
$ws.Range("N80").Value = -8488.143
$ws.Range("H83").Value = 41670496
$ws.Range("J83").Value = 6492.143
$ws.Range("L83").Value = 32460.715
$ws.Range("N83").Value = -42444.715
$ws.Range("H107").Value = 1951.8182
$ws.Range("I107").Value = 574.25
$ws.Range("J107").Value = 2739
$ws.Range("K107").Value = 574.25
$ws.Range("L107").Value = 2739
$ws.Range("M107").Value = 1345.75
$ws.Range("N107").Value = -6579
$ws.Range("H122").Value = 5533.8623
$ws.Range("I122").Value = 6544.65
$ws.Range("J122").Value = 3287.6667
$ws.Range("K122").Value = 19633.95
$ws.Range("L122").Value = 9863.000100000001
$ws.Range("M122").Value = -17183.95
$ws.Range("N122").Value = -14763.0001
$ws.Range("H123").Value = 34999
$ws.Range("J123").Value = 34999
$ws.Range("L123").Value = 34999
$ws.Range("N123").Value = -39899

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3836.5
$ws.Range("I7").Value = 3568.6924
$ws.Range("J7").Value = 4997
$ws.Range("K7").Value = 3568.6924
$ws.Range("L7").Value = 4997
$ws.Range("M7").Value = -3456.6924
$ws.Range("N7").Value = -5221
$ws.Range("H22").Value = 1437.4286
$ws.Range("I22").Value = 638.5
$ws.Range("J22").Value = 2036.625
$ws.Range("K22").Value = 638.5
$ws.Range("L22").Value = 2036.625
$ws.Range("M22").Value = -343.5
$ws.Range("N22").Value = -2626.625
$ws.Range("H27").Value = 1437.4286
$ws.Range("I27").Value = 638.5
$ws.Range("J27").Value = 2036.625
$ws.Range("K27").Value = 638.5
$ws.Range("L27").Value = 2036.625
$ws.Range("M27").Value = -531.5
$ws.Range("N27").Value = -2250.625
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 20000
$ws.Range("K29").Value = 20000
$ws.Range("M29").Value = -19705
$ws.Range("H46").Value = 1140.3096
$ws.Range("I46").Value = 680.1
$ws.Range("J46").Value = 2290.8333
$ws.Range("K46").Value = 680.1
$ws.Range("L46").Value = 2290.8333
$ws.Range("M46").Value = -492.1
$ws.Range("N46").Value = -2666.8333
$ws.Range("H55").Value = 2055.25
$ws.Range("I55").Value = 861.5
$ws.Range("K55").Value = 861.5
$ws.Range("M55").Value = -688.5
$ws.Range("H82").Value = 50001776
$ws.Range("I82").Value = 692.4286
$ws.Range("J82").Value = 166670960
$ws.Range("K82").Value = 692.4286
$ws.Range("L82").Value = 166670960
$ws.Range("M82").Value = -331.4286
$ws.Range("N82").Value = -166671682
$ws.Range("H85").Value = 50001776
$ws.Range("I85").Value = 692.4286
$ws.Range("J85").Value = 166670960
$ws.Range("K85").Value = 692.4286
$ws.Range("L85").Value = 166670960
$ws.Range("M85").Value = 555.5714
$ws.Range("N85").Value = -166673456
$ws.Range("H87").Value = 333358400
$ws.Range("J87").Value = 500025000
$ws.Range("L87").Value = 500025000
$ws.Range("N87").Value = -500027246
$ws.Range("H90").Value = 333358400
$ws.Range("J90").Value = 500025000
$ws.Range("L90").Value = 1500075000
$ws.Range("N90").Value = -1500086232
$ws.Range("H100").Value = 55562536
$ws.Range("I100").Value = 41673012
$ws.Range("K100").Value = 41673012
$ws.Range("M100").Value = -41672471
$ws.Range("H126").Value = 3836.5
$ws.Range("I126").Value = 3568.6924
$ws.Range("J126").Value = 4997
$ws.Range("K126").Value = 10706.0772
$ws.Range("L126").Value = 14991
$ws.Range("M126").Value = -8236.0772
$ws.Range("N126").Value = -19931
$ws.Range("H132").Value = 1100.3334
$ws.Range("I132").Value = 1100.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3301.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -771.0001999999999
$ws.Range("N132").ClearContents()

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 7
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 7
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 7
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -291
$ws.Range("H55").Value = 31166
$ws.Range("I55").Value = 3500
$ws.Range("J55").Value = 44999
$ws.Range("K55").Value = 3500
$ws.Range("L55").Value = 44999
$ws.Range("M55").Value = -3223
$ws.Range("N55").Value = -45553
$ws.Range("H94").Value = 78429.8
$ws.Range("J94").Value = 78429.8
$ws.Range("L94").Value = 78429.8
$ws.Range("N94").Value = -80231.8
$ws.Range("H132").Value = 14503525
$ws.Range("I132").Value = 15162435
$ws.Range("K132").Value = 45487305
$ws.Range("M132").Value = -45484775
$ws.Range("H133").Value = 37274.57
$ws.Range("J133").Value = 37274.57
$ws.Range("L133").Value = 37274.57
$ws.Range("N133").Value = -47394.57
